$d = $word.ActiveDocument
$xml = @'
<w:p w14:paraId="77402162" w14:textId="77777777" w:rsidR="00443E2A" w:rsidRPr="00443E2A" w:rsidRDefault="00443E2A" w:rsidP="00443E2A"/>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Resource: Notas de Estudio - Resúmenes de Introducción de Libros (Tyndale)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>License Information</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Notas de Estudio - Resúmenes de Introducción de Libros (Tyndale)</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> (Spanish) is based on</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Tyndale Open Study Notes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:hyperlink r:id="rId13">
        <w:r>
          <w:rPr>
            <w:lang w:val="es_ES" w:bidi="es_ES"/>
          </w:rPr>
          <w:rPr>
            <w:color w:val="0000EE"/>
            <w:u w:val="single"/>
          </w:rPr>
          <w:t>Tyndale House Publishers</w:t>
        </w:r>
      </w:hyperlink>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t xml:space="preserve">, 2019, which is licensed under a </w:t>
      </w:r>
      <w:hyperlink r:id="rId14">
        <w:r>
          <w:rPr>
            <w:lang w:val="es_ES" w:bidi="es_ES"/>
          </w:rPr>
          <w:rPr>
            <w:color w:val="0000EE"/>
            <w:u w:val="single"/>
          </w:rPr>
          <w:t>CC BY-SA 4.0 license</w:t>
        </w:r>
      </w:hyperlink>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>This PDF version is provided under the same license.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Notas de Estudio - Resúmenes de Introducción de Libros (Tyndale)</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:sectPr>
          <w:headerReference w:type="even" r:id="rId8"/>
          <w:headerReference w:type="default" r:id="rId9"/>
          <w:footerReference w:type="default" r:id="rId10"/>
          <w:footnotePr>
            <w:numRestart w:val="eachSect"/>
          </w:footnotePr>
          <w:type w:val="continuous"/>
          <w:pgSz w:w="12240" w:h="15840" w:code="1"/>
          <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="504" w:footer="504" w:gutter="0"/>
          <w:cols w:space="708"/>
          <w:titlePg/>
          <w:docGrid w:linePitch="360"/>
        </w:sectPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>NUM</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Números</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:sectPr>
          <w:footnotePr>
            <w:numRestart w:val="eachSect"/>
          </w:footnotePr>
          <w:type w:val="continuous"/>
          <w:pgSz w:w="12240" w:h="15840" w:code="1"/>
          <w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="504" w:footer="504" w:gutter="0"/>
          <w:cols w:space="708"/>
          <w:titlePg/>
          <w:docGrid w:linePitch="360"/>
        </w:sectPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Números</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>El Libro de los Números</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Propósito</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Documentar la transición de la generación infiel que salió de Egipto a la generación que entraría en la tierra prometida</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Autor</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Moisés, según la tradición</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Fecha</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Registra eventos que probablemente ocurrieron alrededor de 1445–1406 a.C. o 1269–1230 a.C.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>Contexto</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr/>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
        <w:t>El libro comienza con los israelitas en el Sinaí, sigue su viaje por el desierto y termina con ellos en las llanuras de Moab al borde de la tierra prometida</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es_ES" w:bidi="es_ES"/>
        </w:rPr>
      </w:r>
    </w:p>
    
'@
$d.Content.InsertXML($xml)
$d.Paragraphs.Item(3).Range.ParagraphFormat.Style = "Heading2"
$d.Paragraphs.Item(4).Range.ParagraphFormat.Style = "Heading2"
$d.Paragraphs.Item(7).Range.ParagraphFormat.Style = "Heading1"
$d.Paragraphs.Item(10).Range.ParagraphFormat.Style = "Heading2"
$d.Paragraphs.Item(15).Range.ParagraphFormat.Style = "Heading2"
$d.Paragraphs.Item(17).Range.ParagraphFormat.Style = "Heading3"
$d.Paragraphs.Item(19).Range.ParagraphFormat.Style = "Heading3"
$d.Paragraphs.Item(21).Range.ParagraphFormat.Style = "Heading3"
$d.Paragraphs.Item(23).Range.ParagraphFormat.Style = "Heading3"
Write-Host "Replaced body content and restyled headings"
